# Insert a new weekly price record as row 509 in the daily "Brócoli" data
# sheet. Excel shifts all subsequent rows (old 509..567) down by one
# (new rows 510..568), which matches the target diff exactly (each row's
# values become what the row above used to hold), and the sheet's
# dimension grows from A1:R567 to A1:R568.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 509..567 down to 510..568, creating a blank row 509.
$ws.Rows.Item(509).Insert()

# Populate the new row 509 with the new record's data.
$ws.Range("A509").Value2 = 7
$ws.Range("B509").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C509").Value2 = "Ñuble"
$ws.Range("D509").Value2 = 45142
$ws.Range("D509").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E509").Value2 = 16
$ws.Range("F509").Value2 = 100112023
$ws.Range("G509").Value2 = "Brócoli"
$ws.Range("H509").Value2 = "Sin especificar"
$ws.Range("I509").Value2 = "Segunda"
$ws.Range("J509").Value2 = 500
$ws.Range("K509").Value2 = 800
$ws.Range("L509").Value2 = 800
$ws.Range("M509").Value2 = 800
$ws.Range("N509").Value2 = "$/unidad"
$ws.Range("O509").Value2 = "Provincia de Diguillín"
$ws.Range("P509").Value2 = 800
$ws.Range("Q509").Value2 = 1
$ws.Range("R509").Value2 = "Hortaliza"
